# Apply the weekly data update:
# Insert a new row at row 25 (pushing existing rows 25-92 down to 26-93)
# and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 25; existing row 25 (and all below) shift down by one.
$ws.Rows("25:25").Insert()

# Fill in the new row 25 with the new record's data.
$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "Macroferia Regional de Talca"
$ws.Range("C25").Value = "Maule"
$ws.Range("D25").Value = 44525
$ws.Range("E25").Value = 7
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100108
$ws.Range("H25").Value = "Tropicales y subtropicales"
$ws.Range("I25").Value = 100108002
$ws.Range("J25").Value = "Mango"
$ws.Range("K25").Value = "Sin especificar"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 120
$ws.Range("N25").Value = 6000
$ws.Range("O25").Value = 6000
$ws.Range("P25").Value = 6000
$ws.Range("Q25").Value = "`$/bandeja 4 kilos"
$ws.Range("R25").Value = "Perú"
$ws.Range("S25").Value = 1500
$ws.Range("T25").Value = 4
